{"js": "// Add short gameplay blurbs after each room-name list item, per the\n// \"Updated room descriptions - please add stuff!\" commit.\n//\n// Most rooms simply get a \" - <blurb>\" suffix appended to the existing\n// run. Two rooms also get their label text itself corrected/renamed:\n//   \"Cemetary\"       -> \"Cemetery - bury victims\"\n//   \"Bad guy's home\" -> \"Boss MOB's home - try to defeat him here (if he's there)\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Suffix to append (via InsertLocation.End) to the room-name paragraph.\nconst appendMap = {\n  \"Casino\": \" - sit back and relax, make/lose money\",\n  \"Mansion\": \" - break in, search for items/money/victims\",\n  \"Kitchen\": \" - eat food, look for items/weapons such as knives\",\n  \"Stables\": \" - \",\n  \"Bedroom\": \" - sleep (gain energy)\",\n  \"Hidden room\": \" - find hidden items/powerups - should be difficult to get into(some sort of challenge)\",\n  \"Dining Room\": \" - eat\",\n  \"Jail Room\": \" - go to jail - lose the game here?\",\n  \"Basement\": \" - store victims, weapons, etc.\",\n  \"Attic\": \" - \",\n  \"Pool\": \" - \",\n  \"Train\": \" - travel to a faraway place\",\n  \"Plane\": \" - travel to a faraway place\",\n  \"Paris\": \" - faraway place 1\",\n  \"Spain\": \" - faraway place 2\",\n  \"Dubai\": \" - faraway place 3\",\n  \"Woods\": \" - do creepy things here - idk\",\n  \"Submarine\": \" - travel discreetly \",\n  \"Abandoned Factory\": \" - more creepy stuff here \",\n  \"Court room\": \" - uh oh, someone's in trouble - step 1 of losing\",\n  \"Dark Alleyway\": \" - bad things happen here - item transfers, etc.\",\n  \"Subway station\": \" - \",\n  \"Coffee shop\": \" - more food/energy/caffeine\",\n  \"Hotel\": \" - sleep = more energy\",\n  \"Theme Park\": \" - have fun? find victims\",\n  \"Library\": \" - great spot for a chasing scene - lots of shelves and hiding spots\",\n  \"Diner\": \" - food\",\n  \"Desert\": \" - avoid this room - will drain your energy and water/food levels\",\n};\n\n// Whole-paragraph-text replacements (label text changes + blurb in one go).\nconst replaceMap = {\n  \"Cemetary\": \"Cemetery - bury victims\",\n  \"Bad guy's home\": \"Boss MOB's home - try to defeat him here (if he's there)\",\n};\n\nfor (const p of paragraphs.items) {\n  const current = p.text;\n  if (Object.prototype.hasOwnProperty.call(replaceMap, current)) {\n    p.insertText(replaceMap[current], Word.InsertLocation.replace);\n  } else if (Object.prototype.hasOwnProperty.call(appendMap, current)) {\n    p.insertText(appendMap[current], Word.InsertLocation.end);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Add short gameplay blurbs after each room-name list item, per the\n# \"Updated room descriptions - please add stuff!\" commit.\n#\n# Most rooms simply get a \" - <blurb>\" suffix appended to the existing\n# paragraph. Two rooms also get their label text itself corrected/renamed:\n#   \"Cemetary\"       -> \"Cemetery - bury victims\"\n#   \"Bad guy's home\" -> \"Boss MOB's home - try to defeat him here (if he's there)\"\n\n$d = $word.ActiveDocument\n\n# Suffix to append to the end of the room-name paragraph.\n$appendMap = @{\n    \"Casino\"             = \" - sit back and relax, make/lose money\"\n    \"Mansion\"            = \" - break in, search for items/money/victims\"\n    \"Kitchen\"            = \" - eat food, look for items/weapons such as knives\"\n    \"Stables\"            = \" - \"\n    \"Bedroom\"            = \" - sleep (gain energy)\"\n    \"Hidden room\"        = \" - find hidden items/powerups - should be difficult to get into(some sort of challenge)\"\n    \"Dining Room\"        = \" - eat\"\n    \"Jail Room\"          = \" - go to jail - lose the game here?\"\n    \"Basement\"           = \" - store victims, weapons, etc.\"\n    \"Attic\"              = \" - \"\n    \"Pool\"               = \" - \"\n    \"Train\"              = \" - travel to a faraway place\"\n    \"Plane\"              = \" - travel to a faraway place\"\n    \"Paris\"              = \" - faraway place 1\"\n    \"Spain\"              = \" - faraway place 2\"\n    \"Dubai\"              = \" - faraway place 3\"\n    \"Woods\"              = \" - do creepy things here - idk\"\n    \"Submarine\"          = \" - travel discreetly \"\n    \"Abandoned Factory\"  = \" - more creepy stuff here \"\n    \"Court room\"         = \" - uh oh, someone's in trouble - step 1 of losing\"\n    \"Dark Alleyway\"      = \" - bad things happen here - item transfers, etc.\"\n    \"Subway station\"     = \" - \"\n    \"Coffee shop\"        = \" - more food/energy/caffeine\"\n    \"Hotel\"              = \" - sleep = more energy\"\n    \"Theme Park\"         = \" - have fun? find victims\"\n    \"Library\"            = \" - great spot for a chasing scene - lots of shelves and hiding spots\"\n    \"Diner\"              = \" - food\"\n    \"Desert\"             = \" - avoid this room - will drain your energy and water/food levels\"\n}\n\n# Whole-paragraph-text replacements (label text changes + blurb in one go).\n$replaceMap = @{\n    \"Cemetary\"       = \"Cemetery - bury victims\"\n    \"Bad guy's home\" = \"Boss MOB's home - try to defeat him here (if he's there)\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $current = $p.Range.Text.TrimEnd([char]13)\n    if ($replaceMap.ContainsKey($current)) {\n        $p.Range.Text = $replaceMap[$current]\n    } elseif ($appendMap.ContainsKey($current)) {\n        $p.Range.InsertAfter($appendMap[$current])\n    }\n}\n"}
